# Implement DSA / tornado plot cost values in the 8D cost utilisation
# intervention workbook: update the "cost" column (F) values for the
# relevant rows across all five intervention blocks, and clear the
# explicit cell-level style so the cells fall back to the column's
# default formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CostCell {
    param(
        [string]$Address,
        [double]$NewValue
    )
    $cell = $ws.Range($Address)
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

# Block 1 (rows 2-12)
Set-CostCell "F2"  858.82
Set-CostCell "F4"  128.1
Set-CostCell "F7"  155.94
Set-CostCell "F8"  751.9
Set-CostCell "F9"  751.9
Set-CostCell "F10" 751.9
Set-CostCell "F11" 1938.76

# Block 2 (rows 14-24)
Set-CostCell "F14" 858.82
Set-CostCell "F16" 128.1
Set-CostCell "F19" 155.94
Set-CostCell "F20" 751.9
Set-CostCell "F21" 751.9
Set-CostCell "F22" 751.9
Set-CostCell "F23" 1938.76

# Block 3 (rows 26-36)
Set-CostCell "F26" 858.82
Set-CostCell "F28" 128.1
Set-CostCell "F31" 155.94
Set-CostCell "F32" 751.9
Set-CostCell "F33" 751.9
Set-CostCell "F34" 751.9
Set-CostCell "F35" 1938.76

# Block 4 (rows 38-48)
Set-CostCell "F38" 858.82
Set-CostCell "F40" 128.1
Set-CostCell "F43" 155.94
Set-CostCell "F44" 751.9
Set-CostCell "F45" 751.9
Set-CostCell "F46" 751.9
Set-CostCell "F47" 1938.76

# Block 5 (rows 50-60)
Set-CostCell "F50" 858.82
Set-CostCell "F52" 128.1
Set-CostCell "F55" 155.94
Set-CostCell "F56" 751.9
Set-CostCell "F57" 751.9
Set-CostCell "F58" 751.9
Set-CostCell "F59" 1938.76

# Update the sheet's scroll position / active selection to match the
# state the workbook was left in after the edit.
$ws.Range("J61").Select()
